$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# Update Date (B8) on the Metadata sheet
$wsMeta.Range("B8").Value = "2026-01-14T15:34:52+00:00"

# Update Description (B12) on the Metadata sheet - drop the leading
# "Entrée " so it now reads the same as the Elements sheet's existing
# "Statut clinique du patient" string (which becomes a duplicate and is
# removed from the shared-string table, along with the cells that
# reference it now pointing at the single shared value)
$wsMeta.Range("B12").Value = "Statut clinique du patient"
$wsElements.Range("M2").Value = "Statut clinique du patient"
$wsElements.Range("L6").Value = "Statut clinique du patient"
$wsElements.Range("M6").Value = "Statut clinique du patient"
